# eventAction.xlsx edit: add "taskList" related event rows (259-265) and
# fix the main-menu column-B width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reference cell that already carries the "Chinese text" style (s="1",
# fontId 3 / 宋体) used throughout column B/C of this sheet. We copy its
# format onto the new cells that need the same style instead of creating a
# brand-new font entry (which `.Font.Name = ...` would do).
$styleSource = $ws.Range("B5")

function Set-StyledText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $text
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

function Set-StyledBlank($cellRef) {
    $cell = $ws.Range($cellRef)
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

# Row 259 : taskList
$ws.Range("A259").Value = "taskList"
Set-StyledText "B259" "城市任务"
Set-StyledText "C259" "eventList"
$ws.Range("D259").Value = "closeWindow;taskListCondition"

# Row 260 : taskListCondition
$ws.Range("A260").Value = "taskListCondition"
Set-StyledBlank "B260"
Set-StyledText "C260" "condition"
$ws.Range("D260").Value = "cityHasTasks;cityHasTasksEvent;cityHasNoTasksEvent"

# Row 261 : cityHasTasksEvent
$ws.Range("A261").Value = "cityHasTasksEvent"
Set-StyledText "C261" "eventList"
$ws.Range("D261").Value = "cityHasTasksDialog;cityTasksShowUp"

# Row 262 : cityTasksShowUp
$ws.Range("A262").Value = "cityTasksShowUp"
Set-StyledText "B262" "显示城市任务;特殊选择"
Set-StyledText "C262" "cityTask"
$ws.Range("D262").Value = "shop"

# Row 263 : cityHasTasksDialog
$ws.Range("A263").Value = "cityHasTasksDialog"
Set-StyledText "C263" "dialog"
$ws.Range("D263").Value = "dialog_city_has_tasks"

# Row 264 : cityHasNoTasksEvent
$ws.Range("A264").Value = "cityHasNoTasksEvent"
Set-StyledText "C264" "eventList"
$ws.Range("D264").Value = "cityHasNoTasksDialog;"

# Row 265 : cityHasNoTasksDialog
$ws.Range("A265").Value = "cityHasNoTasksDialog"
Set-StyledText "C265" "dialog"
$ws.Range("D265").Value = "dialog_city_has_no_tasks;shop"

# --- Fix main menu bug: column B used to be a fixed 50.5-char width; now it
# is sized to best-fit the (now shorter) longest entry.
$ws.Columns.Item(2).ColumnWidth = 36.33

# --- Restore the user's view/selection to the tail of the sheet, matching
# where they were working.
$ws.Range("D263").Select() | Out-Null

Write-Output "edit applied"
